$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Runmode changes from N to Y (Description text stays the same)
$ws.Range("C2").Value = "Y"

# Row 3: "My Account" / "Account details" -> "ProductListPage" / "PLP details"
$ws.Range("A3").Value = "ProductListPage"
$ws.Range("B3").Value = "PLP details"
$ws.Range("C3").Value = "Y"

# Row 4 ("Check out" / "CheckoutScenarios" / N) is removed entirely
$ws.Rows.Item(4).Delete()

# Update selection to match target state
$ws.Range("B3").Select()
